# Scheduled-runner market-data refresh for the per-job Leve profit sheets.
# Updates currentAveragePrice(NQ/HQ) + dependent Leve price/profit columns
# (H:N) for the Leve rows whose crafted-item market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4255.4443
$ws.Cells.Item(74, 9).Value = 4000.3333
$ws.Cells.Item(74, 10).Value = 4383
$ws.Cells.Item(74, 11).Value = 4000.3333
$ws.Cells.Item(74, 12).Value = 4383
$ws.Cells.Item(74, 13).Value = -3064.3333
$ws.Cells.Item(74, 14).Value = -6255
$ws.Cells.Item(76, 8).Value = 3115.7896
$ws.Cells.Item(76, 9).Value = 3105.8823
$ws.Cells.Item(76, 11).Value = 3105.8823
$ws.Cells.Item(76, 13).Value = -2790.8823
$ws.Cells.Item(77, 8).Value = 4255.4443
$ws.Cells.Item(77, 9).Value = 4000.3333
$ws.Cells.Item(77, 10).Value = 4383
$ws.Cells.Item(77, 11).Value = 20001.6665
$ws.Cells.Item(77, 12).Value = 21915
$ws.Cells.Item(77, 13).Value = -15321.6665
$ws.Cells.Item(77, 14).Value = -31275
$ws.Cells.Item(79, 8).Value = 3115.7896
$ws.Cells.Item(79, 9).Value = 3105.8823
$ws.Cells.Item(79, 11).Value = 3105.8823
$ws.Cells.Item(79, 13).Value = -2013.8823
$ws.Cells.Item(80, 8).Value = 4535.577
$ws.Cells.Item(80, 9).Value = 3500.6365
$ws.Cells.Item(80, 10).Value = 5294.533
$ws.Cells.Item(80, 11).Value = 10501.9095
$ws.Cells.Item(80, 12).Value = 15883.599
$ws.Cells.Item(80, 13).Value = -9503.9095
$ws.Cells.Item(80, 14).Value = -17879.599
$ws.Cells.Item(83, 8).Value = 4535.577
$ws.Cells.Item(83, 9).Value = 3500.6365
$ws.Cells.Item(83, 10).Value = 5294.533
$ws.Cells.Item(83, 11).Value = 31505.7285
$ws.Cells.Item(83, 12).Value = 47650.79700000001
$ws.Cells.Item(83, 13).Value = -26513.7285
$ws.Cells.Item(83, 14).Value = -57634.79700000001
$ws.Cells.Item(106, 8).Value = 4873.65
$ws.Cells.Item(106, 9).Value = 3995
$ws.Cells.Item(106, 10).Value = 5592.5454
$ws.Cells.Item(106, 11).Value = 3995
$ws.Cells.Item(106, 12).Value = 5592.5454
$ws.Cells.Item(106, 13).Value = -3364
$ws.Cells.Item(106, 14).Value = -6854.5454
$ws.Cells.Item(125, 8).Value = 1305.1428
$ws.Cells.Item(125, 9).Value = 1227.3334
$ws.Cells.Item(125, 10).Value = 1326.3636
$ws.Cells.Item(125, 11).Value = 11046.0006
$ws.Cells.Item(125, 12).Value = 11937.2724
$ws.Cells.Item(125, 13).Value = -8586.000599999999
$ws.Cells.Item(125, 14).Value = -16857.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39, 8).Value = 5008
$ws.Cells.Item(39, 9).Value = 5008
$ws.Cells.Item(39, 11).Value = 5008
$ws.Cells.Item(39, 13).Value = -4488
$ws.Cells.Item(61, 8).Value = 1153.2727
$ws.Cells.Item(61, 9).Value = 1142.1111
$ws.Cells.Item(61, 10).Value = 1203.5
$ws.Cells.Item(61, 11).Value = 1142.1111
$ws.Cells.Item(61, 12).Value = 1203.5
$ws.Cells.Item(61, 13).Value = -930.1111000000001
$ws.Cells.Item(61, 14).Value = -1627.5
$ws.Cells.Item(74, 8).Value = 60277.176
$ws.Cells.Item(74, 9).Value = 84492.664
$ws.Cells.Item(74, 10).Value = 2160
$ws.Cells.Item(74, 11).Value = 84492.664
$ws.Cells.Item(74, 12).Value = 2160
$ws.Cells.Item(74, 13).Value = -83618.664
$ws.Cells.Item(74, 14).Value = -3908
$ws.Cells.Item(77, 8).Value = 60277.176
$ws.Cells.Item(77, 9).Value = 84492.664
$ws.Cells.Item(77, 10).Value = 2160
$ws.Cells.Item(77, 11).Value = 422463.32
$ws.Cells.Item(77, 12).Value = 10800
$ws.Cells.Item(77, 13).Value = -418095.32
$ws.Cells.Item(77, 14).Value = -19536
$ws.Cells.Item(88, 8).Value = 2353.238
$ws.Cells.Item(88, 10).Value = 2541.6667
$ws.Cells.Item(88, 12).Value = 2541.6667
$ws.Cells.Item(88, 14).Value = -3353.6667
$ws.Cells.Item(91, 8).Value = 2353.238
$ws.Cells.Item(91, 10).Value = 2541.6667
$ws.Cells.Item(91, 12).Value = 2541.6667
$ws.Cells.Item(91, 14).Value = -5349.6667
$ws.Cells.Item(122, 8).Value = 2346.0908
$ws.Cells.Item(122, 9).Value = 1843.5555
$ws.Cells.Item(122, 10).Value = 2694
$ws.Cells.Item(122, 11).Value = 5530.666499999999
$ws.Cells.Item(122, 12).Value = 8082
$ws.Cells.Item(122, 13).Value = -3080.666499999999
$ws.Cells.Item(122, 14).Value = -12982
$ws.Cells.Item(136, 8).Value = 1153.2727
$ws.Cells.Item(136, 9).Value = 1142.1111
$ws.Cells.Item(136, 10).Value = 1203.5
$ws.Cells.Item(136, 11).Value = 3426.3333
$ws.Cells.Item(136, 12).Value = 3610.5
$ws.Cells.Item(136, 13).Value = -876.3333000000002
$ws.Cells.Item(136, 14).Value = -8710.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1427.9445
$ws.Cells.Item(99, 9).Value = 1093.5555
$ws.Cells.Item(99, 10).Value = 2431.111
$ws.Cells.Item(99, 11).Value = 1093.5555
$ws.Cells.Item(99, 12).Value = 2431.111
$ws.Cells.Item(99, 13).Value = 404.4445000000001
$ws.Cells.Item(99, 14).Value = -5427.111
$ws.Cells.Item(105, 8).Value = 2149.6
$ws.Cells.Item(105, 10).Value = 2475.2942
$ws.Cells.Item(105, 12).Value = 2475.2942
$ws.Cells.Item(105, 14).Value = -5969.2942

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 3328.1428
$ws.Cells.Item(4, 10).Value = 3328.1428
$ws.Cells.Item(4, 12).Value = 3328.1428
$ws.Cells.Item(4, 14).Value = -3552.1428
$ws.Cells.Item(7, 8).Value = 68.5
$ws.Cells.Item(7, 9).Value = 29.5
$ws.Cells.Item(7, 10).Value = 107.5
$ws.Cells.Item(7, 11).Value = 29.5
$ws.Cells.Item(7, 12).Value = 107.5
$ws.Cells.Item(7, 13).Value = 83.5
$ws.Cells.Item(7, 14).Value = -333.5
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()
$ws.Cells.Item(105, 8).Value = 2680
$ws.Cells.Item(105, 9).Value = 2680
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 2680
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -933
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2053.3333
$ws.Cells.Item(132, 9).Value = 1901.9354
$ws.Cells.Item(132, 10).Value = 2388.5715
$ws.Cells.Item(132, 11).Value = 5705.8062
$ws.Cells.Item(132, 12).Value = 7165.7145
$ws.Cells.Item(132, 13).Value = -3175.8062
$ws.Cells.Item(132, 14).Value = -12225.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 3278.5715
$ws.Cells.Item(51, 9).Value = 400
$ws.Cells.Item(51, 11).Value = 1200
$ws.Cells.Item(51, 13).Value = -740
$ws.Cells.Item(100, 8).Value = 3682.5
$ws.Cells.Item(100, 10).Value = 3682.5
$ws.Cells.Item(100, 12).Value = 11047.5
$ws.Cells.Item(100, 14).Value = -12669.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(112, 8).Value = 30000
$ws.Cells.Item(112, 10).Value = 30000
$ws.Cells.Item(112, 12).Value = 30000
$ws.Cells.Item(112, 14).Value = -32216
$ws.Cells.Item(113, 8).Value = 1615.0435
$ws.Cells.Item(113, 9).Value = 1229.091
$ws.Cells.Item(113, 10).Value = 1968.8334
$ws.Cells.Item(113, 11).Value = 1229.091
$ws.Cells.Item(113, 12).Value = 1968.8334
$ws.Cells.Item(113, 13).Value = 940.9090000000001
$ws.Cells.Item(113, 14).Value = -6308.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(59, 8).Value = 28296
$ws.Cells.Item(59, 10).Value = 28296
$ws.Cells.Item(59, 12).Value = 28296
$ws.Cells.Item(59, 14).Value = -29604
$ws.Cells.Item(68, 8).Value = 7799.4736
$ws.Cells.Item(68, 9).Value = 9737.691999999999
$ws.Cells.Item(68, 10).Value = 3600
$ws.Cells.Item(68, 11).Value = 9737.691999999999
$ws.Cells.Item(68, 12).Value = 3600
$ws.Cells.Item(68, 13).Value = -8988.691999999999
$ws.Cells.Item(68, 14).Value = -5098
$ws.Cells.Item(71, 8).Value = 7799.4736
$ws.Cells.Item(71, 9).Value = 9737.691999999999
$ws.Cells.Item(71, 10).Value = 3600
$ws.Cells.Item(71, 11).Value = 48688.45999999999
$ws.Cells.Item(71, 12).Value = 18000
$ws.Cells.Item(71, 13).Value = -44944.45999999999
$ws.Cells.Item(71, 14).Value = -25488
$ws.Cells.Item(132, 8).Value = 2771.1482
$ws.Cells.Item(132, 9).Value = 2058.2856
$ws.Cells.Item(132, 10).Value = 3020.65
$ws.Cells.Item(132, 11).Value = 6174.8568
$ws.Cells.Item(132, 12).Value = 9061.950000000001
$ws.Cells.Item(132, 13).Value = -3644.8568
$ws.Cells.Item(132, 14).Value = -14121.95

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 550.2308
$ws.Cells.Item(113, 9).Value = 426.66666
$ws.Cells.Item(113, 10).Value = 828.25
$ws.Cells.Item(113, 11).Value = 1279.99998
$ws.Cells.Item(113, 12).Value = 2484.75
$ws.Cells.Item(113, 13).Value = 890.0000199999999
$ws.Cells.Item(113, 14).Value = -6824.75
$ws.Cells.Item(114, 8).Value = 30000
$ws.Cells.Item(114, 10).Value = 30000
$ws.Cells.Item(114, 12).Value = 30000
$ws.Cells.Item(114, 14).Value = -38678
